$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B16 was stored as an inline string "5"; convert it to a real number.
$ws.Range("B16").Value = 5

# Append new row 17 with the new annotation data.
$ws.Range("A17").Value = "Ying Tang"
# B17's "3" must remain text (not become a number) - use a leading
# apostrophe to force text entry, then clear the resulting "quote
# prefix" style so the cell keeps the sheet's default (unstyled) look.
$ws.Range("B17").Value = "'3"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "无"
$ws.Range("D17").Value = "CRT"
$ws.Range("E17").Value = "WRI"
$ws.Range("F17").Value = "4cbdf296-0ef7-4a60-9d08-bf70fb941ab3"
$ws.Range("G17").Value = "SJTB5GZCb_annotated.xlsx"
$ws.Range("H17").Value = "The paper does not sufficiently discuss and compare the relevant neuroscience literature and related work."
